$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.170.04'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.056.48'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.13'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.667'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.44'
$ws.Range('E7').Value = '  +8.88%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.21'
$ws.Range('E12').Value = '  +7.01%  '
$ws.Range('D13').Value = '2.355.59'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.822'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.71'
$ws.Range('E15').Value = '  +8.92%  '
$ws.Range('D16').Value = '2.055.19'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.01'
$ws.Range('E17').Value = '  +26.89%  '
$ws.Range('D18').Value = '37.166.53'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '75.21'
$ws.Range('E19').Value = '  +3.95%  '
$ws.Range('D20').Value = '0.0₃0903'
$ws.Range('E20').Value = '  -3.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.43'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '239.67'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.21'
$ws.Range('E25').Value = '  +11.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.90'
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('E27').Value = '  +3.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.01'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('E29').Value = '  +1.83%  '
$ws.Range('E30').Value = '  +8.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.82'
$ws.Range('E31').Value = '  +5.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0619'
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.54'
$ws.Range('E33').Value = '  +3.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0900'
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.73'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.110'
$ws.Range('E38').Value = '  +5.27%  '
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('E40').Value = '  +31.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.17'
$ws.Range('E41').Value = '  +13.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.03'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0223'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.73'
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.14'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.48'
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('D47').Value = '1.294.59'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.84'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').Value = '2.241.17'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.52'
$ws.Range('E51').Value = '  -20.01%  '
